# Insert a new column before column N ("Late") on the "Repayment schedule" sheet,
# matching the formatting/width of the preceding column (M), leaving the new
# column's cells blank. Then activate this sheet and select cell R6, matching
# where the user ended up after performing the edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

$ws.Activate()
$ws.Range("R6").Select()
